$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "lanka"
$ws.Range("B1").Value = "satish"
$ws.Range("C1").Value = "jjjj"
$ws.Range("C1").Select()
